$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated values for B2:F6 (PME, PMR, PMP, Ciclo Operacional, Ciclo Financeiro)
# Row 2 - 2020
$ws.Range("B2").Value = 295.0753466857058
$ws.Range("C2").Value = 22.43136051933171
$ws.Range("D2").Value = 39.2081063328914
$ws.Range("E2").Value = 317.5067072050375
$ws.Range("F2").Value = 278.2986008721461

# Row 3 - 2021
$ws.Range("B3").Value = 221.9139782122953
$ws.Range("C3").Value = 18.72789121413091
$ws.Range("D3").Value = 67.7835802119571
$ws.Range("E3").Value = 240.6418694264262
$ws.Range("F3").Value = 172.8582892144692

# Row 4 - 2022
$ws.Range("B4").Value = 246.4804909563539
$ws.Range("C4").Value = 21.46119840758994
$ws.Range("D4").Value = 79.28005899172722
$ws.Range("E4").Value = 267.9416893639438
$ws.Range("F4").Value = 188.6616303722166

# Row 5 - 2023
$ws.Range("B5").Value = 214.7684719989315
$ws.Range("C5").Value = 15.23606688231796
$ws.Range("D5").Value = 67.03043013044287
$ws.Range("E5").Value = 230.0045388812495
$ws.Range("F5").Value = 162.9741087508066

# Row 6 - 2024
$ws.Range("B6").Value = 228.5011046094028
$ws.Range("C6").Value = 21.24509240497746
$ws.Range("D6").Value = 67.36648728855037
$ws.Range("E6").Value = 249.7461970143803
$ws.Range("F6").Value = 182.3797097258299
